# ADD results from server
# Update computed values on the "2025", "2030" and "2035" sheets (row 2)
# to reflect newly supplied server results.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 46067.92104640001
$ws.Range("E2").Value = 254981.2308191619
$ws.Range("G2").Value = 64767.40570129279
$ws.Range("I2").Value = 168410.8627624734
$ws.Range("L2").Value = 389668.043028528
$ws.Range("N2").Value = 58533.57433750998
$ws.Range("O2").Value = 57675.15229854788

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 74461.72699460713
$ws.Range("E2").Value = 237095.1674128622
$ws.Range("I2").Value = 200896.4414953778
$ws.Range("L2").Value = 213430.223024568
$ws.Range("N2").Value = 26186.29843170313
$ws.Range("O2").Value = 10012.69360089747

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 9052.525562708204
$ws.Range("B2").Value = 29037.59611842986
$ws.Range("E2").Value = 129622.1712440273
$ws.Range("I2").Value = 134411.1596091401
$ws.Range("N2").Value = 39962.63554520423
$ws.Range("O2").Value = 44812.74248186876
